$d = $word.ActiveDocument

# Find the paragraph containing "ASIGNATURA:" and append " FÍSICA III" as a
# new run right after the existing "ASIGNATURA:" run, inside the same
# paragraph, matching the formatting of the existing run (Verdana, bold, 20).
$found = $d.Content.Find.Execute("ASIGNATURA:", $true, $false, $false, $false,
                                  $false, $true, 1, $false, "ASIGNATURA: FÍSICA III", 2)
